# The commit removes every slide from the deck (slide1.xml .. slide9.xml are
# all deleted, and <p:sldIdLst> disappears from ppt/presentation.xml along
# with it once the slide list is empty). Reproduce that with the PowerPoint
# object model by deleting every slide from the active presentation.

$p = $ppt.ActivePresentation

Write-Output ("Slides before: " + $p.Slides.Count)

# Walk backwards so removing a slide never invalidates the index of the
# slide(s) still to be processed.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $p.Slides.Item($i).Delete()
}

Write-Output ("Slides after: " + $p.Slides.Count)
